$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1295766007657061
$ws.Range("H2").Value = 54.50112707202018
$ws.Range("I2").Value = -19.90039311705729
$ws.Range("G3").Value = 0.08654585096470094
$ws.Range("H3").Value = -25.61959938488688
$ws.Range("G4").Value = -0.01514577055578493
$ws.Range("H4").Value = 1.812013998235429
$ws.Range("G5").Value = -0.01055607088495581
$ws.Range("H5").Value = 85.30852478071137
$ws.Range("G6").Value = 0.0102587413743002
$ws.Range("H6").Value = -71.1392616846786
$ws.Range("G7").Value = -0.01844930161503396
$ws.Range("H7").Value = -190.9052093462718
$ws.Range("G8").Value = -0.1285538300994089
$ws.Range("H8").Value = 8.890328583607266
$ws.Range("G9").Value = -0.1578205217684469
$ws.Range("H9").Value = -15.29254785890846
$ws.Range("G10").Value = -0.08622343307281512
$ws.Range("H10").Value = 20.18353301597519
$ws.Range("G11").Value = -0.09963493035276311
$ws.Range("H11").Value = -49.73495668883047
$ws.Range("G12").Value = -0.3428913147196478
$ws.Range("H12").Value = 17.41465533710561
$ws.Range("G13").Value = -0.4179487337133465
$ws.Range("H13").Value = 6.894786481650642
$ws.Range("G14").Value = -0.04262981171984808
$ws.Range("H14").Value = 15.97668886031958
$ws.Range("G15").Value = 0.010575109525961
$ws.Range("H15").Value = 112.7886391353961
$ws.Range("G16").Value = 0.1226285622704251
$ws.Range("H16").Value = -16.02088291793515
$ws.Range("G17").Value = 0.1111008878412645
$ws.Range("H17").Value = -9.395251913725776
$ws.Range("G18").Value = 0.1473807317572547
$ws.Range("H18").Value = 6.444989758877581
$ws.Range("G19").Value = 0.1324612766226856
$ws.Range("H19").Value = 38.83848242368666
$ws.Range("G20").Value = 0.04549785822051627
$ws.Range("H20").Value = 77.48991834845495
$ws.Range("G21").Value = 0.06738738007493135
$ws.Range("H21").Value = -10.34991526230334
$ws.Range("G24").Value = 0.07912888497988632
$ws.Range("H24").Value = -21.21501696425947
$ws.Range("G25").Value = 0.1737047210155188
$ws.Range("H25").Value = 14.62137498416663
$ws.Range("G26").Value = 0.08933131092336503
$ws.Range("H26").Value = 12.89938281290031
$ws.Range("G27").Value = 0.07132575000409781
$ws.Range("H27").Value = -28.5689757215738
$ws.Range("G28").Value = -0.2573594666503971
$ws.Range("H28").Value = -20.75847003798419
$ws.Range("G29").Value = -0.2338037909518537
$ws.Range("H29").Value = -13.88026715952888
$ws.Range("G30").Value = 0.07063171359994788
$ws.Range("H30").Value = 60.0477783006878
$ws.Range("G31").Value = 0.03297166118520949
$ws.Range("H31").Value = 25.20989304137859
$ws.Range("G32").Value = 0.097275197777203
$ws.Range("H32").Value = 2.437528024332915
$ws.Range("G33").Value = 0.1263673864480193
$ws.Range("H33").Value = 21.54814009649536
$ws.Range("G34").Value = 0.009951934715990065
$ws.Range("H34").Value = -78.56465734798552
$ws.Range("G35").Value = 0.01078005874187858
$ws.Range("H35").Value = 42.28231439995581
$ws.Range("G36").Value = 0.04612902812932372
$ws.Range("H36").Value = -20.10386649301405
$ws.Range("G37").Value = 0.06105204567560074
$ws.Range("H37").Value = -13.18728873846296
$ws.Range("G38").Value = 0.02943902727003059
$ws.Range("H38").Value = -43.80407723869907
$ws.Range("G39").Value = -0.006968889235226058
$ws.Range("H39").Value = -133.6055290709452
$ws.Range("G40").Value = -0.001498239553577904
$ws.Range("H40").Value = 82.35494756051803
$ws.Range("G41").Value = 0.01033712871504971
$ws.Range("H41").Value = -70.7626105939442
$ws.Range("G42").Value = 0.143523481285987
$ws.Range("H42").Value = 7.358546328619962
$ws.Range("G43").Value = 0.1416298539810153
$ws.Range("H43").Value = -4.931307066158563
$ws.Range("G44").Value = 0.004278316780011757
$ws.Range("H44").Value = 150.2652434244137
$ws.Range("G45").Value = -0.01443886878337824
$ws.Range("H45").Value = -31.50315331535647
$ws.Range("G46").Value = -0.007487949896170274
$ws.Range("H46").Value = -127.3867447873703
$ws.Range("G47").Value = -0.02997099399722834
$ws.Range("H47").Value = -223.0161670814284
$ws.Range("G48").Value = 0.03631718205628851
$ws.Range("H48").Value = -27.75878990861587
$ws.Range("G49").Value = 0.05966364319624113
$ws.Range("H49").Value = -9.689187014783474
$ws.Range("G50").Value = 0.1477215052010647
$ws.Range("H50").Value = -8.394127393011752
$ws.Range("G51").Value = 0.172421096467501
$ws.Range("H51").Value = 0.7633096485598335
$ws.Range("G52").Value = -0.1686432639670065
$ws.Range("H52").Value = -5.118874947397918
$ws.Range("G53").Value = -0.1317493171948181
$ws.Range("H53").Value = -4.517741789952728
$ws.Range("G54").Value = 0.09058584982275487
$ws.Range("H54").Value = -3.346633821976666
$ws.Range("G55").Value = 0.1089275019187739
$ws.Range("H55").Value = -3.672813319448658
$ws.Range("G56").Value = 0.004812471240225346
$ws.Range("H56").Value = 165.9172970603223
$ws.Range("G57").Value = -0.03797900105718639
$ws.Range("H57").Value = -66.08366322420586
$ws.Range("G58").Value = 0.04343918477082789
$ws.Range("H58").Value = -22.96374067631849
$ws.Range("G59").Value = 0.06198573640409088
$ws.Range("H59").Value = -13.69284280543179
$ws.Range("G60").Value = 0.06911924485010899
$ws.Range("H60").Value = -1.22146677759126
$ws.Range("G61").Value = 0.05555100145432924
$ws.Range("H61").Value = 16.88356541391079
$ws.Range("G62").Value = 0.05800430210312708
$ws.Range("H62").Value = -20.49884631656348
$ws.Range("G63").Value = 0.06075397396511949
$ws.Range("H63").Value = -7.094788098448153
$ws.Range("G64").Value = -0.03076785143736567
$ws.Range("H64").Value = 25.70934596055264
$ws.Range("G65").Value = -0.00724295378970915
$ws.Range("H65").Value = 85.31547175007083
$ws.Range("G66").Value = 0.04269997272669277
$ws.Range("H66").Value = 125.5203770365715
$ws.Range("G67").Value = 0.0284804414586097
$ws.Range("H67").Value = 8.912297426866916
$ws.Range("G68").Value = -0.01210071706313794
$ws.Range("H68").Value = -2223.106479811778
$ws.Range("G69").Value = -0.005801458495430863
$ws.Range("H69").Value = 55.07031711170851
$ws.Range("G70").Value = -0.05202482494236257
$ws.Range("H70").Value = -89.53911335798223
$ws.Range("G71").Value = -0.0648519845702232
$ws.Range("H71").Value = -17.71589752237616
$ws.Range("G72").Value = -0.1502832689040412
$ws.Range("H72").Value = -1.314008120500677
$ws.Range("G73").Value = -0.1390408027452464
$ws.Range("H73").Value = 3.970605545665582
$ws.Range("G74").Value = 0.131317032404888
$ws.Range("H74").Value = 4.177612518742702
$ws.Range("G75").Value = 0.1557779351866707
$ws.Range("H75").Value = 15.24425283663594
$ws.Range("G76").Value = -0.05050397108185651
$ws.Range("H76").Value = -46.65210430413831
$ws.Range("G77").Value = -0.06871238910363191
$ws.Range("H77").Value = -48.76995137264375
$ws.Range("G78").Value = 0.08092177855064224
$ws.Range("H78").Value = -12.20524236503486
$ws.Range("G79").Value = 0.09686152802431476
$ws.Range("H79").Value = 0.3700333276133108
$ws.Range("G80").Value = -0.1572107916526153
$ws.Range("H80").Value = 3.202989462354445
$ws.Range("G81").Value = -0.1896877460649929
$ws.Range("H81").Value = 12.3635493633186
$ws.Range("G82").Value = 0.1614688025380438
$ws.Range("H82").Value = 16.37703266802605
$ws.Range("G83").Value = 0.2027966317101158
$ws.Range("H83").Value = 23.19106467018088
$ws.Range("G84").Value = 0.08727317799695559
$ws.Range("H84").Value = 523.6013299892869
$ws.Range("G85").Value = 0.08356361445977097
$ws.Range("H85").Value = 269.1105454638557
